# Aktualisierter Zeitplan (vh)
# Markiert mehrere Aufgaben als zu 100% erledigt und traegt das tatsaechliche
# Fertigstellungsdatum nach.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgabenliste Projekt 1")

# Zeile 26 - "(S) Kamerafahrt-Szene: Bett +Junge" (Tobias): auf 100% setzen.
# Die Fortschritt-Spalte (F) verliert dabei ihre Formel und wird zu einem
# festen Wert (wie im Original-Workbook geschehen).
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1

# Zeile 27 - "(S)Kamerafahrt-Szene: Kommode + Auto + Wecker" (Tobias): auf 100%.
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1

# Zeile 32 - "(A) - Szene 1 - 4 des Kamerafahrt-Videos" (Viktoria): auf 100%
# und tatsaechliches Fertigstellungsdatum (27.11.2015) eintragen.
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 1
$ws.Range("H27").Copy() | Out-Null
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("H32").Value = 42335

# Zeile 33 - "(A) - Szene 5 - 9 des Kamerafahrt-Videos" (Jana): auf 100%,
# Formel in Fortschritt bleibt erhalten.
$ws.Range("E33").Value = 1
$ws.Range("H27").Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null
$ws.Range("H33").Value = 42335

# Zeile 34 - "(A) - Szene 10 - 14 des Kamerafahrt-Videos" (Tobias): auf 100%.
$ws.Range("E34").Value = 1
$ws.Range("H27").Copy() | Out-Null
$ws.Range("H34").PasteSpecial(-4122) | Out-Null
$ws.Range("H34").Value = 42335

# Zeile 35 - "(R) - Rendern des Kamerafahrt-Videos" (Alle): auf 100%.
$ws.Range("E35").Value = 1
$ws.Range("H27").Copy() | Out-Null
$ws.Range("H35").PasteSpecial(-4122) | Out-Null
$ws.Range("H35").Value = 42335

$excel.CutCopyMode = 0

# Ansicht: wieder an den Tabellenanfang scrollen, Auswahl auf H35 (letzte
# bearbeitete Zelle).
$ws.Activate()
$ws.Range("H35").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 1
